# Danh_Muc_Khach_Hang_Import.xlsx
# "fix lay nhat ky mau hoa don va file excel"
#
# Changes applied:
#   1. Rename the single worksheet "KETOAN-BACHKHOA" -> "HOADON-BACHKHOA"
#      (this also updates <sheet name="..."/> in workbook.xml).
#   2. Move the sheet's active selection from N10 to I24 and drop the
#      previous topLeftCell="I1" scroll-freeze (selecting a cell clears the
#      stale top-left anchor the same way real Excel does on click).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet.
$ws.Name = "HOADON-BACHKHOA"

# 2) Update the active cell / selection on the sheet.
$ws.Activate()
$ws.Range("I24").Select()

# Best-effort: the author's diff also shows the cached window geometry
# (bookViews/workbookView windowWidth/windowHeight) growing from
# 23040x10452 to 28800x12144, and the x15ac:absPath "recent folder" hint
# changing from G:\downloads\ to C:\Users\vanlt\Downloads\. Both are
# session/UI metadata that Excel stamps from the OS window manager and the
# file-open dialog rather than values exposed on the Workbook/Window COM
# object model, but we still try the documented Window properties here in
# case the host maps them through - failures are swallowed so the rest of
# the script keeps running.
try {
    $win = $excel.ActiveWindow
    $win.WindowState = -4143  # xlNormal
    $win.Width = 28800
    $win.Height = 12144
} catch {
}
